$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency table: column D is Price, column E is Volume(1h).
# Every cell in D/E is stored as text in the source sheet. Plain strings such as
# "27.251.98" (two dots) are never mistaken for numbers, but some new prices
# (e.g. "210.66") parse as valid numbers, so Excel would silently convert them.
# Prefixing those with an apostrophe forces text entry, then ClearFormats() drops
# the "quote prefix" cell style Excel adds so no formatting otherwise changes.

$ws.Range("D2").Value = "27.251.98"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "1.564.54"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'210.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("D11").Value = "'0.0871"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("D12").Value = "1.787.96"
$ws.Range("D13").Value = "1.576.92"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "'0.518"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "27.196.39"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "'61.90"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "'218.08"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'7.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'9.37"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("D25").Value = "'151.62"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D28").Value = "'15.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "1.456.39"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "'3.18"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  +4.53%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("D40").Value = "'5.88"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").Value = "'0.814"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").Value = "  +1.11%  "
$ws.Range("D44").Value = "'0.984"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").Value = "'64.42"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "1.699.60"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").Value = "'85.89"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("D51").Value = "'0.0945"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.77%  "
